$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# "12/05/2020" looks like a date, so a plain Range.Value assignment would be
# auto-converted to a date serial (and bump the cell's style/number format).
# Route it through a formula-literal + paste-values round trip instead, which
# keeps it a plain text/shared-string value with the cell's original style.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""12/05/2020"""

$scratch.Copy()
$ws.Range("H2").PasteSpecial(-4163)
$ws.Range("H4").PasteSpecial(-4163)
$ws.Range("H6").PasteSpecial(-4163)

$scratch.ClearContents()
